$wb = $excel.ActiveWorkbook

# District Summary
$ws = $wb.Worksheets.Item("District Summary")
$ws.Range("A1").Value = "Total Schools"
$ws.Range("B1").Value = "Total Students"
$ws.Range("C1").Value = "Total Budget"
$ws.Range("D1").Value = "Average Math Score"
$ws.Range("E1").Value = "Average Reading Score"
$ws.Range("F1").Value = "% Passing Math"
$ws.Range("G1").Value = "% Passing Reading"
$ws.Range("H1").Value = "% Overall Passing"

# School Summary
$ws = $wb.Worksheets.Item("School Summary")
$ws.Range("C1").Value = "Total Budget"
$ws.Range("D1").Value = "Average Math Score"
$ws.Range("E1").Value = "Average Reading Score"
$ws.Range("F1").Value = "% Passing Math"
$ws.Range("G1").Value = "% Passing Reading"
$ws.Range("H1").Value = "% Overall Passing"

# Top Performing Schools
$ws = $wb.Worksheets.Item("Top Performing Schools")
$ws.Range("C1").Value = "Total Budget"
$ws.Range("D1").Value = "Average Math Score"
$ws.Range("E1").Value = "Average Reading Score"
$ws.Range("F1").Value = "% Passing Math"
$ws.Range("G1").Value = "% Passing Reading"
$ws.Range("H1").Value = "% Overall Passing"

# Bottom Performing Schools
$ws = $wb.Worksheets.Item("Bottom Performing Schools")
$ws.Range("C1").Value = "Total Budget"
$ws.Range("D1").Value = "Average Math Score"
$ws.Range("E1").Value = "Average Reading Score"
$ws.Range("F1").Value = "% Passing Math"
$ws.Range("G1").Value = "% Passing Reading"
$ws.Range("H1").Value = "% Overall Passing"

# Scores by School Spending
$ws = $wb.Worksheets.Item("Scores by School Spending")
$ws.Range("B1").Value = "Average Math Score"
$ws.Range("C1").Value = "Average Reading Score"
$ws.Range("D1").Value = "% Passing Math"
$ws.Range("E1").Value = "% Passing Reading"
$ws.Range("F1").Value = "% Overall Passing"

# Scores by School Size
$ws = $wb.Worksheets.Item("Scores by School Size")
$ws.Range("B1").Value = "Average Math Score"
$ws.Range("C1").Value = "Average Reading Score"
$ws.Range("D1").Value = "% Passing Math"
$ws.Range("E1").Value = "% Passing Reading"
$ws.Range("F1").Value = "% Overall Passing"

# Scores by School Type
$ws = $wb.Worksheets.Item("Scores by School Type")
$ws.Range("B1").Value = "Average Math Score"
$ws.Range("C1").Value = "Average Reading Score"
$ws.Range("D1").Value = "% Passing Math"
$ws.Range("E1").Value = "% Passing Reading"
$ws.Range("F1").Value = "% Overall Passing"
